# Applies the Dec-20-2022 symbol-list refresh: updated Price (D) values,
# the Worstin24h/Bestin24h badge-text tweaks on two coin labels (E18, E47),
# and flips the Hora flag (G) from 0 to 1 for every data row (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.74"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1"
$ws.Range("G2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.41"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1"
$ws.Range("G3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.300"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "1"
$ws.Range("G4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05560"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "1"
$ws.Range("G5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.376"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "1"
$ws.Range("G6").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.354"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "1"
$ws.Range("G7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8106"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "1"
$ws.Range("G8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9516"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "1"
$ws.Range("G9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1390"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "1"
$ws.Range("G10").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07398"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "1"
$ws.Range("G11").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03125"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "1"
$ws.Range("G12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03038"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "1"
$ws.Range("G13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09299"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1"
$ws.Range("G14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.566"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "1"
$ws.Range("G15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001612"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "1"
$ws.Range("G16").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04699"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "1"
$ws.Range("G17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005740"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "1"
$ws.Range("G18").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006401"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "1"
$ws.Range("G19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005040"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "1"
$ws.Range("G20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001038"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "1"
$ws.Range("G21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001498"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "1"
$ws.Range("G22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.804"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "1"
$ws.Range("G23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.125"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "1"
$ws.Range("G24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3244"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "1"
$ws.Range("G25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1284"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "1"
$ws.Range("G26").NumberFormat = "General"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "1"
$ws.Range("G27").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003094"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "1"
$ws.Range("G28").NumberFormat = "General"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "1"
$ws.Range("G29").NumberFormat = "General"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "1"
$ws.Range("G30").NumberFormat = "General"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "1"
$ws.Range("G31").NumberFormat = "General"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "1"
$ws.Range("G32").NumberFormat = "General"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "1"
$ws.Range("G33").NumberFormat = "General"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "1"
$ws.Range("G34").NumberFormat = "General"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "1"
$ws.Range("G35").NumberFormat = "General"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "1"
$ws.Range("G36").NumberFormat = "General"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "1"
$ws.Range("G37").NumberFormat = "General"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "1"
$ws.Range("G38").NumberFormat = "General"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "1"
$ws.Range("G39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03876"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "1"
$ws.Range("G40").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006900"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "1"
$ws.Range("G41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1043"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "1"
$ws.Range("G42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003020"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "1"
$ws.Range("G43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007700"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "1"
$ws.Range("G44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005807"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "1"
$ws.Range("G45").NumberFormat = "General"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "1"
$ws.Range("G46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005490"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "1"
$ws.Range("G47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6788"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "1"
$ws.Range("G48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1189"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "1"
$ws.Range("G49").NumberFormat = "General"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "1"
$ws.Range("G50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01008"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "1"
$ws.Range("G51").NumberFormat = "General"
